$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after row 11 (for the new "tot_oecd" and "tot_imf" entries),
# shifting the "Liquidity indicators" / "Global conditions" blocks down by 2 rows.
$ws.Rows("12:13").Insert()

# Row 11: rename "tot" -> "tot_eiu"
$ws.Cells.Item(11, 1).Value = "tot_eiu"

# Row 12 (new): tot_oecd
$ws.Cells.Item(12, 1).Value = "tot_oecd"
$ws.Cells.Item(12, 4).Value = "same, but retrieved from OECD database, since data from economist intelligence unit is indexed weirdly (index years are different) - OECD omits malta and cyprus though, but shouldn't matter too much"

# Row 13 (new): tot_imf
$ws.Cells.Item(13, 1).Value = "tot_imf"
$ws.Cells.Item(13, 2).Value = "terms of trade"
$ws.Cells.Item(12, 2).Value = "terms of trade (malta, cyprus missing)"
$ws.Cells.Item(13, 4).Value = "same, own calculation (see terms_of_trade_IMF.xlsx), export/import price indices according to IMF weighted by ratio of exports/imports to GDP - should be fine, I guess"

# Update row 11's remaining description text
$ws.Cells.Item(11, 2).Value = "terms of trade (not the same index across countries?)"
$ws.Cells.Item(11, 4).Value = "increase in terms of trade means that export prices in relation to average import prices increase - decrease default risk since it favors the economy and makes it easier to collect funds for debt servicing"

# Update selection to reflect where the author was last working
$ws.Range("D13").Select()
